$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TIPO_DOCUMENTO_ids_map")
$ws.Copy($null, $ws)

$ws2004 = $wb.Worksheets.Item("TIPO_DOCUMENTO_ids_map")
$ws2004.Name = "TIPO_DOCUMENTO_ids_map_2004"

$ws2017 = $wb.Worksheets.Item("TIPO_DOCUMENTO_ids_map (2)")
$ws2017.Name = "TIPO_DOCUMENTO_ids_map_2017"

$ws2017.Range("A2:A26").ClearContents()
$ws2017.Range("B2:B26").ClearContents()

$ws2017.Activate()
$ws2017.Range("F10").Select()

Write-Output ($wb.Worksheets | ForEach-Object { $_.Name })
